$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 338.6742503333333
$ws.Range("H2").Value = 1016.022751
$ws.Range("I2").Value = 0.5849329800180821
$ws.Range("J2").Value = 0.584932980018082
$ws.Range("M2").Value = 71.44418333333333
$ws.Range("N2").Value = 214.33255
$ws.Range("O2").Value = 0.6986063918429039
$ws.Range("P2").Value = 0.6986063918429037
$ws.Range("Q2").Value = 24196.30523109389
$ws.Range("R2").Value = 217766.7470798451
$ws.Range("S2").Value = 0.4086379186403498
$ws.Range("T2").Value = 0.4086379186403496

$ws.Range("G3").Value = 338.6742503333333
$ws.Range("H3").Value = 1016.022751
$ws.Range("I3").Value = 0.5849329800180821
$ws.Range("J3").Value = 0.584932980018082
$ws.Range("O3").Value = 0.1188372961583501
$ws.Range("P3").Value = 0.1188372961583501
$ws.Range("Q3").Value = 4115.942144617449
$ws.Range("R3").Value = 37043.47930155705
$ws.Range("S3").Value = 0.0695118537791951
$ws.Range("T3").Value = 0.06951185377919507

$ws.Range("G4").Value = 338.6742503333333
$ws.Range("H4").Value = 1016.022751
$ws.Range("I4").Value = 0.5849329800180821
$ws.Range("J4").Value = 0.584932980018082
$ws.Range("M4").Value = 18.571964
$ws.Range("N4").Value = 55.715892
$ws.Range("O4").Value = 0.1816032062252276
$ws.Range("P4").Value = 0.1816032062252276
$ws.Range("Q4").Value = 6289.845984917653
$ws.Range("R4").Value = 56608.61386425888
$ws.Range("S4").Value = 0.1062257045981607
$ws.Range("T4").Value = 0.1062257045981607

$ws.Range("G5").Value = 338.6742503333333
$ws.Range("H5").Value = 1016.022751
$ws.Range("I5").Value = 0.5849329800180821
$ws.Range("J5").Value = 0.584932980018082
$ws.Range("K5").Value = 2.0
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09747100000000002
$ws.Range("N5").Value = 0.292413
$ws.Range("O5").Value = 0.000953105773518577
$ws.Range("P5").Value = 0.0009531057735185768
$ws.Range("Q5").Value = 33.01091785424033
$ws.Range("R5").Value = 297.098260688163
$ws.Range("S5").Value = 0.0005575030003766605
$ws.Range("T5").Value = 0.0005575030003766603

$ws.Range("I6").Value = 0.279688040971731
$ws.Range("J6").Value = 0.2796880409717309
$ws.Range("M6").Value = 71.44418333333333
$ws.Range("N6").Value = 214.33255
$ws.Range("O6").Value = 0.6986063918429039
$ws.Range("P6").Value = 0.6986063918429037
$ws.Range("Q6").Value = 11569.56000092437
$ws.Range("R6").Value = 104126.0400083194
$ws.Range("S6").Value = 0.1953918531448713
$ws.Range("T6").Value = 0.1953918531448712

$ws.Range("I7").Value = 0.279688040971731
$ws.Range("J7").Value = 0.2796880409717309
$ws.Range("O7").Value = 0.1188372961583501
$ws.Range("P7").Value = 0.1188372961583501
$ws.Range("S7").Value = 0.03323737055690635
$ws.Range("T7").Value = 0.03323737055690634

$ws.Range("I8").Value = 0.279688040971731
$ws.Range("J8").Value = 0.2796880409717309
$ws.Range("M8").Value = 18.571964
$ws.Range("N8").Value = 55.715892
$ws.Range("O8").Value = 0.1816032062252276
$ws.Range("P8").Value = 0.1816032062252276
$ws.Range("Q8").Value = 3007.514983137289
$ws.Range("R8").Value = 27067.6348482356
$ws.Range("S8").Value = 0.05079224498331917
$ws.Range("T8").Value = 0.05079224498331916

$ws.Range("I9").Value = 0.279688040971731
$ws.Range("J9").Value = 0.2796880409717309
$ws.Range("K9").Value = 2.0
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09747100000000002
$ws.Range("N9").Value = 0.292413
$ws.Range("O9").Value = 0.000953105773518577
$ws.Range("P9").Value = 0.0009531057735185768
$ws.Range("Q9").Value = 15.78430223757567
$ws.Range("R9").Value = 142.058720138181
$ws.Range("S9").Value = 0.0002665722866342572
$ws.Range("T9").Value = 0.000266572286634257

$ws.Range("G10").Value = 77.79536166666666
$ws.Range("H10").Value = 233.386085
$ws.Range("I10").Value = 0.1343623634996766
$ws.Range("J10").Value = 0.1343623634996766
$ws.Range("M10").Value = 71.44418333333333
$ws.Range("N10").Value = 214.33255
$ws.Range("O10").Value = 0.6986063918429039
$ws.Range("P10").Value = 0.6986063918429037
$ws.Range("Q10").Value = 5558.026081396305
$ws.Range("R10").Value = 50022.23473256675
$ws.Range("S10").Value = 0.09386640596399376
$ws.Range("T10").Value = 0.09386640596399373

$ws.Range("G11").Value = 77.79536166666666
$ws.Range("H11").Value = 233.386085
$ws.Range("I11").Value = 0.1343623634996766
$ws.Range("J11").Value = 0.1343623634996766
$ws.Range("O11").Value = 0.1188372961583501
$ws.Range("P11").Value = 0.1188372961583501
$ws.Range("Q11").Value = 945.4548358029538
$ws.Range("R11").Value = 8509.093522226585
$ws.Range("S11").Value = 0.01596725998374696
$ws.Range("T11").Value = 0.01596725998374695

$ws.Range("G12").Value = 77.79536166666666
$ws.Range("H12").Value = 233.386085
$ws.Range("I12").Value = 0.1343623634996766
$ws.Range("J12").Value = 0.1343623634996766
$ws.Range("M12").Value = 18.571964
$ws.Range("N12").Value = 55.715892
$ws.Range("O12").Value = 0.1816032062252276
$ws.Range("P12").Value = 0.1816032062252276
$ws.Range("Q12").Value = 1444.812656240313
$ws.Range("R12").Value = 13003.31390616282
$ws.Range("S12").Value = 0.02440063600754077
$ws.Range("T12").Value = 0.02440063600754076

$ws.Range("G13").Value = 77.79536166666666
$ws.Range("H13").Value = 233.386085
$ws.Range("I13").Value = 0.1343623634996766
$ws.Range("J13").Value = 0.1343623634996766
$ws.Range("K13").Value = 2.0
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.09747100000000002
$ws.Range("N13").Value = 0.292413
$ws.Range("O13").Value = 0.000953105773518577
$ws.Range("P13").Value = 0.0009531057735185768
$ws.Range("Q13").Value = 7.582791697011667
$ws.Range("R13").Value = 68.245125273105
$ws.Range("S13").Value = 0.0001280615443951435
$ws.Range("T13").Value = 0.0001280615443951434

$ws.Range("G14").Value = 0.5886170000000001
$ws.Range("H14").Value = 1.765851
$ws.Range("I14").Value = 0.001016615510510267
$ws.Range("J14").Value = 0.001016615510510266
$ws.Range("M14").Value = 71.44418333333333
$ws.Range("N14").Value = 214.33255
$ws.Range("O14").Value = 0.6986063918429039
$ws.Range("P14").Value = 0.6986063918429037
$ws.Range("Q14").Value = 42.05326086111667
$ws.Range("R14").Value = 378.47934775005
$ws.Range("S14").Value = 0.0007102140936891091
$ws.Range("T14").Value = 0.0007102140936891088

$ws.Range("G15").Value = 0.5886170000000001
$ws.Range("H15").Value = 1.765851
$ws.Range("I15").Value = 0.001016615510510267
$ws.Range("J15").Value = 0.001016615510510266
$ws.Range("O15").Value = 0.1188372961583501
$ws.Range("P15").Value = 0.1188372961583501
$ws.Range("Q15").Value = 7.153521458905667
$ws.Range("R15").Value = 64.38169313015099
$ws.Range("S15").Value = 0.0001208118385016808
$ws.Range("T15").Value = 0.0001208118385016808

$ws.Range("G16").Value = 0.5886170000000001
$ws.Range("H16").Value = 1.765851
$ws.Range("I16").Value = 0.001016615510510267
$ws.Range("J16").Value = 0.001016615510510266
$ws.Range("M16").Value = 18.571964
$ws.Range("N16").Value = 55.715892
$ws.Range("O16").Value = 0.1816032062252276
$ws.Range("P16").Value = 0.1816032062252276
$ws.Range("Q16").Value = 10.931773733788
$ws.Range("R16").Value = 98.385963604092
$ws.Range("S16").Value = 0.000184620636206961
$ws.Range("T16").Value = 0.0001846206362069609

$ws.Range("G17").Value = 0.5886170000000001
$ws.Range("H17").Value = 1.765851
$ws.Range("I17").Value = 0.001016615510510267
$ws.Range("J17").Value = 0.001016615510510266
$ws.Range("K17").Value = 2.0
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.09747100000000002
$ws.Range("N17").Value = 0.292413
$ws.Range("O17").Value = 0.000953105773518577
$ws.Range("P17").Value = 0.0009531057735185768
$ws.Range("Q17").Value = 0.05737308760700002
$ws.Range("R17").Value = 0.5163577884630001
$ws.Range("S17").Value = 0.0000009689421125158708
$ws.Range("T17").Value = 0.0000009689421125158704
